$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as exact text
# (matches source "inlineStr" cells, avoids Excel auto-converting to a
# Number and silently dropping meaningful trailing zeros / grouping dots).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.160.85'
$ws.Range("E2").Value = '  +0.71%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.519.81'
$ws.Range("E3").Value = '  +0.43%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.30'
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.77'
$ws.Range("E6").Value = '  +2.73%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.594'
$ws.Range("E8").Value = '  +3.67%  '
$ws.Range("E9").Value = '  +7.14%  '
$ws.Range("E10").Value = '  -0.55%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.438'
$ws.Range("E11").Value = '  -0.29%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.130.65'
$ws.Range("E12").Value = '  +0.44%  '
$ws.Range("E13").Value = '  +0.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29.12'
$ws.Range("E14").Value = '  +3.35%  '
$ws.Range("E15").Value = '  +1.71%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.154.94'
$ws.Range("E16").Value = '  +0.73%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.564.04'
$ws.Range("E17").Value = '  +2.69%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.34'
$ws.Range("E18").Value = '  -0.10%  '
$ws.Range("E19").Value = '  +1.35%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '395.82'
$ws.Range("E20").Value = '  +1.08%  '
$ws.Range("E21").Value = '  +0.36%  '
$ws.Range("E22").Value = '  -0.50%  '
$ws.Range("E23").Value = '  +1.62%  '
$ws.Range("E24").Value = '  +0.27%  '
$ws.Range("E25").Value = '  -3.66%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.33'
$ws.Range("E26").Value = '  +1.88%  '
$ws.Range("E27").Value = '  +0.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("E28").Value = '  -0.14%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.29'
$ws.Range("E29").Value = '  -2.27%  '
$ws.Range("E30").Value = '  -0.97%  '
$ws.Range("E31").Value = '  -0.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '23.94'
$ws.Range("E32").Value = '  +1.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.36'
$ws.Range("E33").Value = '  -1.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.64'
$ws.Range("E34").Value = '  +1.49%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '163.12'
$ws.Range("E35").Value = '  +1.30%  '
$ws.Range("E36").Value = '  +0.86%  '
$ws.Range("E37").Value = '  +0.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.01'
$ws.Range("E38").Value = '  +5.25%  '
$ws.Range("E39").Value = '  +1.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0747'
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '26.49'
$ws.Range("E42").Value = '  -0.56%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.64'
$ws.Range("E43").Value = '  +3.86%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.811.42'
$ws.Range("E44").Value = '  -0.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.97'
$ws.Range("E45").Value = '  -0.87%  '
$ws.Range("E46").Value = '  -2.66%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '337.24'
$ws.Range("E47").Value = '  -5.02%  '
$ws.Range("E48").Value = '  +1.20%  '
$ws.Range("E49").Value = '  -1.82%  '
$ws.Range("E50").Value = '  +0.35%  '
$ws.Range("E51").Value = '  -0.44%  '
